# Apply crypto symbol-list update (GitHub Actions refresh, Fri Jan 27 11:13:10 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the numeric-looking columns (D: Price, E: Volume, G: Hora)
# so the values are stored as text, matching the source sheet's inline-string cells.
$ws.Range("D2:E25").NumberFormat = "@"
$ws.Range("D27:E28").NumberFormat = "@"
$ws.Range("D39:E46").NumberFormat = "@"
$ws.Range("D48:E49").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Updated coin rows (values refreshed from coinranking.com)
$ws.Range("D2").Value = '305.25'
$ws.Range("E2").Value = '0.40%'
$ws.Range("G2").Value = '11'
$ws.Range("D3").Value = '35.57'
$ws.Range("E3").Value = '-0.20%'
$ws.Range("G3").Value = '11'
$ws.Range("D4").Value = '5.047'
$ws.Range("E4").Value = '-0.69%'
$ws.Range("G4").Value = '11'
$ws.Range("D5").Value = '0.08005'
$ws.Range("E5").Value = '-0.47%'
$ws.Range("G5").Value = '11'
$ws.Range("D6").Value = '1.859'
$ws.Range("E6").Value = '-2.30%'
$ws.Range("G6").Value = '11'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '4.142'
$ws.Range("E7").Value = '-1.03%'
$ws.Range("G7").Value = '11'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '7.779'
$ws.Range("E8").Value = '0.49%'
$ws.Range("G8").Value = '11'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9199'
$ws.Range("E9").Value = '-0.86%'
$ws.Range("G9").Value = '11'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1268'
$ws.Range("E10").Value = '-8.43%'
$ws.Range("G10").Value = '11'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1922'
$ws.Range("E11").Value = '1.67%'
$ws.Range("G11").Value = '11'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09134'
$ws.Range("E12").Value = '0.08%'
$ws.Range("G12").Value = '11'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03453'
$ws.Range("E13").Value = '-4.38%'
$ws.Range("G13").Value = '11'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09853'
$ws.Range("E14").Value = '0.40%'
$ws.Range("G14").Value = '11'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001407'
$ws.Range("E15").Value = '-1.79%'
$ws.Range("G15").Value = '11'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.006229'
$ws.Range("E16").Value = '5.35%'
$ws.Range("G16").Value = '11'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.848'
$ws.Range("E17").Value = '8.33%'
$ws.Range("G17").Value = '11'
$ws.Range("D18").Value = '3.358'
$ws.Range("E18").Value = '12.78%'
$ws.Range("G18").Value = '11'
$ws.Range("D19").Value = '0.3418'
$ws.Range("E19").Value = '-1.05%'
$ws.Range("G19").Value = '11'
$ws.Range("D20").Value = '0.1320'
$ws.Range("E20").Value = '-0.89%'
$ws.Range("G20").Value = '11'
$ws.Range("D21").Value = '5.215'
$ws.Range("E21").Value = '6.36%'
$ws.Range("G21").Value = '11'
$ws.Range("D22").Value = '0.2305'
$ws.Range("E22").Value = '-8.23%'
$ws.Range("G22").Value = '11'
$ws.Range("D23").Value = '0.04435'
$ws.Range("E23").Value = '-0.33%'
$ws.Range("G23").Value = '11'
$ws.Range("E24").Value = '0.92%'
$ws.Range("G24").Value = '11'
$ws.Range("D25").Value = '0.004888'
$ws.Range("E25").Value = '2.23%'
$ws.Range("G25").Value = '11'
$ws.Range("G26").Value = '11'
$ws.Range("E27").Value = '-19.95%'
$ws.Range("G27").Value = '11'
$ws.Range("E28").Value = '41.82%'
$ws.Range("G28").Value = '11'
$ws.Range("G29").Value = '11'
$ws.Range("G30").Value = '11'
$ws.Range("G31").Value = '11'
$ws.Range("G32").Value = '11'
$ws.Range("G33").Value = '11'
$ws.Range("G34").Value = '11'
$ws.Range("G35").Value = '11'
$ws.Range("G36").Value = '11'
$ws.Range("G37").Value = '11'
$ws.Range("G38").Value = '11'
$ws.Range("D39").Value = '0.01937'
$ws.Range("E39").Value = '-0.74%'
$ws.Range("G39").Value = '11'
$ws.Range("D40").Value = '0.05205'
$ws.Range("E40").Value = '6.73%'
$ws.Range("G40").Value = '11'
$ws.Range("D41").Value = '0.007637'
$ws.Range("E41").Value = '-0.23%'
$ws.Range("G41").Value = '11'
$ws.Range("D42").Value = '0.01018'
$ws.Range("E42").Value = '9.72%'
$ws.Range("G42").Value = '11'
$ws.Range("D43").Value = '0.1351'
$ws.Range("E43").Value = '-1.57%'
$ws.Range("G43").Value = '11'
$ws.Range("D44").Value = '0.002163'
$ws.Range("E44").Value = '2.76%'
$ws.Range("G44").Value = '11'
$ws.Range("D45").Value = '0.009629'
$ws.Range("E45").Value = '-15.64%'
$ws.Range("G45").Value = '11'
$ws.Range("D46").Value = '0.00006195'
$ws.Range("E46").Value = '-3.05%'
$ws.Range("G46").Value = '11'
$ws.Range("G47").Value = '11'
$ws.Range("E48").Value = '0.48%'
$ws.Range("G48").Value = '11'
$ws.Range("E49").Value = '39.14%'
$ws.Range("G49").Value = '11'
$ws.Range("G50").Value = '11'
$ws.Range("G51").Value = '11'
